$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'286.54"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'2.32%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'28.62"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'4.16%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'5.072"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'4.99%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.06649"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'3.59%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'7.379"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'4.73%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'3.405"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'2.98%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'1.373"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'5.64%"
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.9392"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'3.73%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.1578"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'2.72%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.06616"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'7.43%"
$ws.Range("E11").Style = "Normal"
$ws.Range("E12").Value = "'1.35%"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.02937"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'0.39%"
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.08978"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'-0.21%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.001584"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'-0.04%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.04500"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'2.20%"
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'0.0006461"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'0.43%"
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'0.006264"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'4.57%"
$ws.Range("E18").Style = "Normal"
$ws.Range("E20").Value = "'0.75%"
$ws.Range("E20").Style = "Normal"
$ws.Range("E21").Value = "'2.26%"
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'0.1297"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'-4.18%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'4.065"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'4.28%"
$ws.Range("E23").Style = "Normal"
$ws.Range("D25").Value = "'0.001183"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'0.82%"
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.004145"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'-3.57%"
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'0.0001250"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'6.13%"
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'0.0001616"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'-2.46%"
$ws.Range("E28").Style = "Normal"
$ws.Range("D40").Value = "'0.04200"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'3.06%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.006732"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'0.98%"
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").Value = "'-11.47%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.002020"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'-2.70%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.01234"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'12.12%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.00005582"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'0.29%"
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Value = "'25.93%"
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.01305"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'-29.38%"
$ws.Range("E47").Style = "Normal"
